$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear row 2 (A2:B2) back to empty/default state
$ws.Range("A2").Value = $null
$ws.Range("B2").Value = $null

# New form responses shifted into rows 3-6
$ws.Range("A3").Value = (Get-Date -Year 2023 -Month 12 -Day 24 -Hour 18 -Minute 52 -Second 44)
$ws.Range("B3").Value = "הקבוצה של: אור, המפקד, רון"

$ws.Range("A4").Value = (Get-Date -Year 2023 -Month 12 -Day 24 -Hour 18 -Minute 53 -Second 1)
$ws.Range("B4").Value = "הקבוצה של: איי, עמרי"

$ws.Range("A5").Value = (Get-Date -Year 2023 -Month 12 -Day 24 -Hour 18 -Minute 53 -Second 9)
$ws.Range("B5").Value = "הקבוצה של: דור, הקשבי"

$ws.Range("A6").Value = (Get-Date -Year 2023 -Month 12 -Day 24 -Hour 18 -Minute 53 -Second 18)
$ws.Range("B6").Value = "הקבוצה של: לא נתמך עי גוגל, טון"

# Append 4 new blank rows at the bottom (rows 103-106)
$lastRow = 102
for ($i = 1; $i -le 4; $i++) {
    $r = $lastRow + $i
    $ws.Range("A$r`:H$r").Value = $null
}

# Update selection to reflect where the user last clicked
$ws.Range("E12").Select()
